$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (Prudential Financial, Inc. / PRU)
$ws.Range("K2").Value = 59.1
$ws.Range("N2").Value = 85.82376350509293

# Row 3 (UnitedHealth Group Incorporated / UNH)
$ws.Range("K3").Value = 56.1
$ws.Range("N3").Value = 85.82376350509293

# Row 4 (MetLife, Inc. / MET)
$ws.Range("D4").Value = 76.56
$ws.Range("F4").Value = 3.57
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 54.7
$ws.Range("N4").Value = 85.82376350509293

# Row 5 (American International Group, I / AIG)
$ws.Range("K5").Value = 51.7
$ws.Range("N5").Value = 85.82376350509293
